$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for the 2022-Q3 summary figures,
#    pushing the existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Insert a blank row at position 2 (existing rows 2-7 shift to 3-8).
$summary.Rows.Item(2).Insert()

# Clear the formatting Excel auto-propagated into the new B2:D2 cells so
# they end up unstyled, matching the rest of the data rows.
$summary.Range("B2:D2").ClearFormats()

# Give A2 the same "index" style used by the column elsewhere (copy format
# only, from A3, then set its own value).
$summary.Cells.Item(3,1).Copy()
$summary.Cells.Item(2,1).PasteSpecial(-4122)
$summary.Cells.Item(2,1).Value = 0

$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 32
$summary.Cells.Item(2,4).Value = 6.14

# The index column (A) is a simple 0-based row counter; the row that used
# to be last (2021-Q1) moved from row 7 to row 8, so its index must bump
# from 5 to 6.
$summary.Cells.Item(8,1).Value = 6

# ---------------------------------------------------------------------
# 2) Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    sheet that is currently "2022-Q2"), holding the quarter's fund
#    holdings detail.
# ---------------------------------------------------------------------
$beforeSheet = $wb.Worksheets.Item(2)
$ns = $wb.Worksheets.Add($beforeSheet)
$ns.Name = "2022-Q3"

# Apply the header/index box style (bold, thin border, centered) used on
# the other quarter sheets by copying it from the "总计" sheet, BEFORE any
# values are written (so no stray auto styles get introduced first).
$summary.Cells.Item(1,2).Copy()
$ns.Range("B1:H1").PasteSpecial(-4122)
$summary.Cells.Item(2,1).Copy()
$ns.Range("A2:A33").PasteSpecial(-4122)

$ns.Cells.Item(1,2).Value = "基金代码"
$ns.Cells.Item(1,3).Value = "基金名称"
$ns.Cells.Item(1,4).Value = "基金规模"
$ns.Cells.Item(1,5).Value = "股票总仓位"
$ns.Cells.Item(1,6).Value = "仓位占比"
$ns.Cells.Item(1,7).Value = "持有市值(亿元)"
$ns.Cells.Item(1,8).Value = "仓位排名"
$ns.Cells.Item(2,1).Value = 0
$ns.Cells.Item(2,2).Value = "'900099"
$ns.Cells.Item(2,3).Value = "中信证券红利价值一年持有混合B"
$ns.Cells.Item(2,4).Value = "'53.18"
$ns.Cells.Item(2,5).Value = "'89.06"
$ns.Cells.Item(2,6).Value = "'2.88"
$ns.Cells.Item(2,7).Value = "'1.5316"
$ns.Cells.Item(2,8).Value = 7

$ns.Cells.Item(3,1).Value = 1
$ns.Cells.Item(3,2).Value = "'013516"
$ns.Cells.Item(3,3).Value = "汇添富蓝筹稳健混合E"
$ns.Cells.Item(3,4).Value = "'62.33"
$ns.Cells.Item(3,5).Value = "'71.66"
$ns.Cells.Item(3,6).Value = "'2.21"
$ns.Cells.Item(3,7).Value = "'1.3775"
$ns.Cells.Item(3,8).Value = 10

$ns.Cells.Item(4,1).Value = 2
$ns.Cells.Item(4,2).Value = "'519066"
$ns.Cells.Item(4,3).Value = "汇添富蓝筹稳健混合A"
$ns.Cells.Item(4,4).Value = "'62.32"
$ns.Cells.Item(4,5).Value = "'71.66"
$ns.Cells.Item(4,6).Value = "'2.21"
$ns.Cells.Item(4,7).Value = "'1.3773"
$ns.Cells.Item(4,8).Value = 10

$ns.Cells.Item(5,1).Value = 3
$ns.Cells.Item(5,2).Value = "'900089"
$ns.Cells.Item(5,3).Value = "中信证券红利价值一年持有混合C"
$ns.Cells.Item(5,4).Value = "'18.22"
$ns.Cells.Item(5,5).Value = "'89.06"
$ns.Cells.Item(5,6).Value = "'2.88"
$ns.Cells.Item(5,7).Value = "'0.5247"
$ns.Cells.Item(5,8).Value = 7

$ns.Cells.Item(6,1).Value = 4
$ns.Cells.Item(6,2).Value = "'000478"
$ns.Cells.Item(6,3).Value = "建信中证500指数增强A"
$ns.Cells.Item(6,4).Value = "'45.95"
$ns.Cells.Item(6,5).Value = "'82.53"
$ns.Cells.Item(6,6).Value = "'0.79"
$ns.Cells.Item(6,7).Value = "'0.3630"
$ns.Cells.Item(6,8).Value = 10

$ns.Cells.Item(7,1).Value = 5
$ns.Cells.Item(7,2).Value = "'004148"
$ns.Cells.Item(7,3).Value = "圆信永丰多策略精选混合"
$ns.Cells.Item(7,4).Value = "'5.72"
$ns.Cells.Item(7,5).Value = "'91.99"
$ns.Cells.Item(7,6).Value = "'4.39"
$ns.Cells.Item(7,7).Value = "'0.2511"
$ns.Cells.Item(7,8).Value = 8

$ns.Cells.Item(8,1).Value = 6
$ns.Cells.Item(8,2).Value = "'001490"
$ns.Cells.Item(8,3).Value = "汇添富国企创新增长股票A"
$ns.Cells.Item(8,4).Value = "'6.95"
$ns.Cells.Item(8,5).Value = "'84.32"
$ns.Cells.Item(8,6).Value = "'2.48"
$ns.Cells.Item(8,7).Value = "'0.1724"
$ns.Cells.Item(8,8).Value = 10

$ns.Cells.Item(9,1).Value = 7
$ns.Cells.Item(9,2).Value = "'515760"
$ns.Cells.Item(9,3).Value = "华夏中证浙江国资创新发展ETF"
$ns.Cells.Item(9,4).Value = "'2.04"
$ns.Cells.Item(9,5).Value = "'99.57"
$ns.Cells.Item(9,6).Value = "'5.90"
$ns.Cells.Item(9,7).Value = "'0.1204"
$ns.Cells.Item(9,8).Value = 5

$ns.Cells.Item(10,1).Value = 8
$ns.Cells.Item(10,2).Value = "'900011"
$ns.Cells.Item(10,3).Value = "中信证券红利价值一年持有混合A"
$ns.Cells.Item(10,4).Value = "'3.05"
$ns.Cells.Item(10,5).Value = "'89.06"
$ns.Cells.Item(10,6).Value = "'2.88"
$ns.Cells.Item(10,7).Value = "'0.0878"
$ns.Cells.Item(10,8).Value = 7

$ns.Cells.Item(11,1).Value = 9
$ns.Cells.Item(11,2).Value = "'016854"
$ns.Cells.Item(11,3).Value = "汇添富成长多因子量化策略股票C"
$ns.Cells.Item(11,4).Value = "'8.11"
$ns.Cells.Item(11,5).Value = "'93.37"
$ns.Cells.Item(11,6).Value = "'1.01"
$ns.Cells.Item(11,7).Value = "'0.0819"
$ns.Cells.Item(11,8).Value = 6

$ns.Cells.Item(12,1).Value = 10
$ns.Cells.Item(12,2).Value = "'005062"
$ns.Cells.Item(12,3).Value = "博时中证500指数增强A"
$ns.Cells.Item(12,4).Value = "'3.36"
$ns.Cells.Item(12,5).Value = "'90.74"
$ns.Cells.Item(12,6).Value = "'1.63"
$ns.Cells.Item(12,7).Value = "'0.0548"
$ns.Cells.Item(12,8).Value = 9

$ns.Cells.Item(13,1).Value = 11
$ns.Cells.Item(13,2).Value = "'005351"
$ns.Cells.Item(13,3).Value = "汇添富行业整合主题混合A"
$ns.Cells.Item(13,4).Value = "'0.67"
$ns.Cells.Item(13,5).Value = "'83.14"
$ns.Cells.Item(13,6).Value = "'5.20"
$ns.Cells.Item(13,7).Value = "'0.0348"
$ns.Cells.Item(13,8).Value = 5

$ns.Cells.Item(14,1).Value = 12
$ns.Cells.Item(14,2).Value = "'006969"
$ns.Cells.Item(14,3).Value = "圆信永丰高端制造混合"
$ns.Cells.Item(14,4).Value = "'0.91"
$ns.Cells.Item(14,5).Value = "'87.79"
$ns.Cells.Item(14,6).Value = "'3.67"
$ns.Cells.Item(14,7).Value = "'0.0334"
$ns.Cells.Item(14,8).Value = 3

$ns.Cells.Item(15,1).Value = 13
$ns.Cells.Item(15,2).Value = "'005633"
$ns.Cells.Item(15,3).Value = "建信中证500指数增强C"
$ns.Cells.Item(15,4).Value = "'3.42"
$ns.Cells.Item(15,5).Value = "'82.53"
$ns.Cells.Item(15,6).Value = "'0.79"
$ns.Cells.Item(15,7).Value = "'0.0270"
$ns.Cells.Item(15,8).Value = 10

$ns.Cells.Item(16,1).Value = 14
$ns.Cells.Item(16,2).Value = "'512190"
$ns.Cells.Item(16,3).Value = "浙商汇金中证浙江凤凰行动50ETF"
$ns.Cells.Item(16,4).Value = "'0.48"
$ns.Cells.Item(16,5).Value = "'98.92"
$ns.Cells.Item(16,6).Value = "'4.98"
$ns.Cells.Item(16,7).Value = "'0.0239"
$ns.Cells.Item(16,8).Value = 4

$ns.Cells.Item(17,1).Value = 15
$ns.Cells.Item(17,2).Value = "'620002"
$ns.Cells.Item(17,3).Value = "金元顺安成长动力混合"
$ns.Cells.Item(17,4).Value = "'0.51"
$ns.Cells.Item(17,5).Value = "'62.21"
$ns.Cells.Item(17,6).Value = "'3.95"
$ns.Cells.Item(17,7).Value = "'0.0201"
$ns.Cells.Item(17,8).Value = 2

$ns.Cells.Item(18,1).Value = 16
$ns.Cells.Item(18,2).Value = "'005795"
$ns.Cells.Item(18,3).Value = "博时中证500指数增强C"
$ns.Cells.Item(18,4).Value = "'1.11"
$ns.Cells.Item(18,5).Value = "'90.74"
$ns.Cells.Item(18,6).Value = "'1.63"
$ns.Cells.Item(18,7).Value = "'0.0181"
$ns.Cells.Item(18,8).Value = 9

$ns.Cells.Item(19,1).Value = 17
$ns.Cells.Item(19,2).Value = "'015453"
$ns.Cells.Item(19,3).Value = "中欧中证500指数增强A"
$ns.Cells.Item(19,4).Value = "'1.20"
$ns.Cells.Item(19,5).Value = "'88.73"
$ns.Cells.Item(19,6).Value = "'1.40"
$ns.Cells.Item(19,7).Value = "'0.0168"
$ns.Cells.Item(19,8).Value = 7

$ns.Cells.Item(20,1).Value = 18
$ns.Cells.Item(20,2).Value = "'006522"
$ns.Cells.Item(20,3).Value = "财通新兴蓝筹混合A"
$ns.Cells.Item(20,4).Value = "'0.22"
$ns.Cells.Item(20,5).Value = "'94.22"
$ns.Cells.Item(20,6).Value = "'3.58"
$ns.Cells.Item(20,7).Value = "'0.0079"
$ns.Cells.Item(20,8).Value = 9

$ns.Cells.Item(21,1).Value = 19
$ns.Cells.Item(21,2).Value = "'015454"
$ns.Cells.Item(21,3).Value = "中欧中证500指数增强C"
$ns.Cells.Item(21,4).Value = "'0.34"
$ns.Cells.Item(21,5).Value = "'88.73"
$ns.Cells.Item(21,6).Value = "'1.40"
$ns.Cells.Item(21,7).Value = "'0.0048"
$ns.Cells.Item(21,8).Value = 7

$ns.Cells.Item(22,1).Value = 20
$ns.Cells.Item(22,2).Value = "'015123"
$ns.Cells.Item(22,3).Value = "汇添富国企创新增长股票C"
$ns.Cells.Item(22,4).Value = "'0.15"
$ns.Cells.Item(22,5).Value = "'84.32"
$ns.Cells.Item(22,6).Value = "'2.48"
$ns.Cells.Item(22,7).Value = "'0.0037"
$ns.Cells.Item(22,8).Value = 10

$ns.Cells.Item(23,1).Value = 21
$ns.Cells.Item(23,2).Value = "'015225"
$ns.Cells.Item(23,3).Value = "汇添富中证细分化工产业主题指数增强A"
$ns.Cells.Item(23,4).Value = "'0.11"
$ns.Cells.Item(23,5).Value = "'92.25"
$ns.Cells.Item(23,6).Value = "'3.11"
$ns.Cells.Item(23,7).Value = "'0.0034"
$ns.Cells.Item(23,8).Value = 7

$ns.Cells.Item(24,1).Value = 22
$ns.Cells.Item(24,2).Value = "'005260"
$ns.Cells.Item(24,3).Value = "银华稳健增利灵活配置混合A"
$ns.Cells.Item(24,4).Value = "'0.28"
$ns.Cells.Item(24,5).Value = "'91.67"
$ns.Cells.Item(24,6).Value = "'0.96"
$ns.Cells.Item(24,7).Value = "'0.0027"
$ns.Cells.Item(24,8).Value = 4

$ns.Cells.Item(25,1).Value = 23
$ns.Cells.Item(25,2).Value = "'515510"
$ns.Cells.Item(25,3).Value = "嘉实中证500成长估值ETF"
$ns.Cells.Item(25,4).Value = "'0.13"
$ns.Cells.Item(25,5).Value = "'98.46"
$ns.Cells.Item(25,6).Value = "'1.21"
$ns.Cells.Item(25,7).Value = "'0.0016"
$ns.Cells.Item(25,8).Value = 10

$ns.Cells.Item(26,1).Value = 24
$ns.Cells.Item(26,2).Value = "'006523"
$ns.Cells.Item(26,3).Value = "财通新兴蓝筹混合C"
$ns.Cells.Item(26,4).Value = "'0.03"
$ns.Cells.Item(26,5).Value = "'94.22"
$ns.Cells.Item(26,6).Value = "'3.58"
$ns.Cells.Item(26,7).Value = "'0.0011"
$ns.Cells.Item(26,8).Value = 9

$ns.Cells.Item(27,1).Value = 25
$ns.Cells.Item(27,2).Value = "'015226"
$ns.Cells.Item(27,3).Value = "汇添富中证细分化工产业主题指数增强C"
$ns.Cells.Item(27,4).Value = "'0.02"
$ns.Cells.Item(27,5).Value = "'92.25"
$ns.Cells.Item(27,6).Value = "'3.11"
$ns.Cells.Item(27,7).Value = "'0.0006"
$ns.Cells.Item(27,8).Value = 7

$ns.Cells.Item(28,1).Value = 26
$ns.Cells.Item(28,2).Value = "'015124"
$ns.Cells.Item(28,3).Value = "汇添富国企创新增长股票D"
$ns.Cells.Item(28,4).Value = "'0.01"
$ns.Cells.Item(28,5).Value = "'84.32"
$ns.Cells.Item(28,6).Value = "'2.48"
$ns.Cells.Item(28,7).Value = "'0.0002"
$ns.Cells.Item(28,8).Value = 10

$ns.Cells.Item(29,1).Value = 27
$ns.Cells.Item(29,2).Value = "'005261"
$ns.Cells.Item(29,3).Value = "银华稳健增利灵活配置混合C"
$ns.Cells.Item(29,4).Value = "'0.02"
$ns.Cells.Item(29,5).Value = "'91.67"
$ns.Cells.Item(29,6).Value = "'0.96"
$ns.Cells.Item(29,7).Value = "'0.0002"
$ns.Cells.Item(29,8).Value = 4

$ns.Cells.Item(30,1).Value = 28
$ns.Cells.Item(30,2).Value = "'001050"
$ns.Cells.Item(30,3).Value = "汇添富成长多因子量化策略股票A"
$ns.Cells.Item(30,4).Value = "'0.00"
$ns.Cells.Item(30,5).Value = "'93.37"
$ns.Cells.Item(30,6).Value = "'1.01"
$ns.Cells.Item(30,7).Value = 0
$ns.Cells.Item(30,8).Value = 6

$ns.Cells.Item(31,1).Value = 29
$ns.Cells.Item(31,2).Value = "'013515"
$ns.Cells.Item(31,3).Value = "汇添富蓝筹稳健混合C"
$ns.Cells.Item(31,4).Value = "'0.00"
$ns.Cells.Item(31,5).Value = "'71.66"
$ns.Cells.Item(31,6).Value = "'2.21"
$ns.Cells.Item(31,7).Value = 0
$ns.Cells.Item(31,8).Value = 10

$ns.Cells.Item(32,1).Value = 30
$ns.Cells.Item(32,2).Value = "'015191"
$ns.Cells.Item(32,3).Value = "汇添富行业整合主题混合D"
$ns.Cells.Item(32,4).Value = "'0.00"
$ns.Cells.Item(32,5).Value = "'83.14"
$ns.Cells.Item(32,6).Value = "'5.20"
$ns.Cells.Item(32,7).Value = 0
$ns.Cells.Item(32,8).Value = 5

$ns.Cells.Item(33,1).Value = 31
$ns.Cells.Item(33,2).Value = "'015190"
$ns.Cells.Item(33,3).Value = "汇添富行业整合主题混合C"
$ns.Cells.Item(33,4).Value = "'0.00"
$ns.Cells.Item(33,5).Value = "'83.14"
$ns.Cells.Item(33,6).Value = "'5.20"
$ns.Cells.Item(33,7).Value = 0
$ns.Cells.Item(33,8).Value = 5

# Remove any stray quote-prefix / number-format styling the data writes
# below may introduce on the text columns (B-G), leaving them with the
# workbook's default (unstyled) cell format.
$ns.Range("B2:G33").ClearFormats()

Write-Host "Edit complete"
